$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate a paragraph whose visible text contains a given substring.
# ---------------------------------------------------------------------------
function Find-ParagraphContaining($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) "Tools: Postman, Chrome DevTools" paragraph - bold/enlarge the paragraph
#    mark (the formatting that trails the final run) from 12pt to 14pt bold,
#    matching the "Methods:" / "11. Screenshots..." heading runs elsewhere
#    in the document. The visible runs themselves keep their own formatting.
# ---------------------------------------------------------------------------
$toolsPara = Find-ParagraphContaining("Tools: Postman, Chrome DevTools")
if ($toolsPara -ne $null) {
    $body = '<w:p w:rsidR="00FD681D" w:rsidRDefault="00FD681D" w:rsidP="00BD3B53">' + `
            '<w:pPr><w:rPr>' + `
            '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
            '<w:b/><w:sz w:val="28"/><w:szCs w:val="28"/>' + `
            '</w:rPr></w:pPr>' + `
            '<w:r w:rsidRPr="00FD681D"><w:rPr>' + `
            '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
            '<w:b/><w:sz w:val="28"/><w:szCs w:val="28"/>' + `
            '</w:rPr><w:t xml:space="preserve">Tools: </w:t></w:r>' + `
            '<w:r w:rsidRPr="00FD681D"><w:rPr>' + `
            '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
            '<w:sz w:val="24"/><w:szCs w:val="24"/>' + `
            '</w:rPr><w:t xml:space="preserve">Postman, Chrome </w:t></w:r>' + `
            '<w:proofErr w:type="spellStart"/>' + `
            '<w:r w:rsidRPr="00FD681D"><w:rPr>' + `
            '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
            '<w:sz w:val="24"/><w:szCs w:val="24"/>' + `
            '</w:rPr><w:t>DevTools</w:t></w:r>' + `
            '<w:proofErr w:type="spellEnd"/>' + `
            '</w:p>'
    $null = $toolsPara.Range.InsertXML($xmlHeader + $body + $xmlFooter)
}

# ---------------------------------------------------------------------------
# 2) Screenshot/demo link paragraph - swap the Google Drive URL for the new
#    one and make the whole line bold/14pt (paragraph mark + run), matching
#    the "11. Screenshots or Demo link:" heading above it.
# ---------------------------------------------------------------------------
$linkPara = Find-ParagraphContaining("drive.google.com")
if ($linkPara -ne $null) {
    $newUrl = "https://drive.google.com/file/d/1pgkPU-QDRdV4TGYKn6iSZaM4ZH9ZLN4O/view?usp=drivesdk"
    $body = '<w:p w:rsidR="00FD681D" w:rsidRPr="00016F5B" w:rsidRDefault="002F2EC7" w:rsidP="00BD3B53">' + `
            '<w:pPr><w:rPr>' + `
            '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
            '<w:b/><w:sz w:val="28"/><w:szCs w:val="28"/>' + `
            '</w:rPr></w:pPr>' + `
            '<w:r w:rsidRPr="00016F5B"><w:rPr>' + `
            '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
            '<w:b/><w:sz w:val="28"/><w:szCs w:val="28"/>' + `
            '</w:rPr><w:t>' + $newUrl + '</w:t></w:r>' + `
            '</w:p>'
    $null = $linkPara.Range.InsertXML($xmlHeader + $body + $xmlFooter)
}

Write-Output "edit complete"
